$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '77.305.34'
$ws.Range("E2").Value = '  +1.23%  '

# Row 3
$ws.Range("D3").Value = '3.141.60'
$ws.Range("E3").Value = '  +5.97%  '

# Row 4
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '201.82'
$ws.Range("E5").Value = '  +0.96%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '627.17'
$ws.Range("E6").Value = '  -0.69%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.00%  '

# Row 8
$ws.Range("E8").Value = '  +14.54%  '

# Row 9
$ws.Range("E9").Value = '  +4.64%  '

# Row 10
$ws.Range("D10").Value = '3.138.96'
$ws.Range("E10").Value = '  +5.91%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.539'
$ws.Range("E11").Value = '  +24.81%  '

# Row 12
$ws.Range("E12").Value = '  +1.22%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.45'
$ws.Range("E13").Value = '  +9.09%  '

# Row 14
$ws.Range("D14").Value = '3.718.49'
$ws.Range("E14").Value = '  +5.82%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000225'
$ws.Range("E15").Value = '  +20.59%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '30.65'
$ws.Range("E16").Value = '  +6.01%  '

# Row 17
$ws.Range("D17").Value = '77.157.28'
$ws.Range("E17").Value = '  +1.16%  '

# Row 18
$ws.Range("D18").Value = '3.137.75'
$ws.Range("E18").Value = '  +5.91%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.03'
$ws.Range("E19").Value = '  +5.04%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.41'
$ws.Range("E20").Value = '  +7.88%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '426.25'
$ws.Range("E21").Value = '  +14.95%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.83'
$ws.Range("E22").Value = '  +25.71%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.81'
$ws.Range("E23").Value = '  +12.45%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.84'
$ws.Range("E24").Value = '  +6.93%  '

# Row 25
$ws.Range("D25").Value = '3.302.49'
$ws.Range("E25").Value = '  +5.76%  '

# Row 26
$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.67'
$ws.Range("E26").Value = '  +9.38%  '

# Row 27
$ws.Range("B27").Value = 'Litecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '75.53'
$ws.Range("E27").Value = '  +3.98%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.73'
$ws.Range("E28").Value = '  +11.49%  '

# Row 29
$ws.Range("E29").Value = '  +0.20%  '

# Row 30
$ws.Range("E30").Value = '  +8.62%  '

# Row 31
$ws.Range("E31").Value = '  +0.11%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.78'
$ws.Range("E32").Value = '  +7.17%  '

# Row 33
$ws.Range("E33").Value = '  +7.51%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '519.59'
$ws.Range("E34").Value = '  +2.65%  '

# Row 35
$ws.Range("E35").Value = '  +1.38%  '

# Row 36
$ws.Range("E36").Value = '  +22.88%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '22.34'
$ws.Range("E37").Value = '  +10.64%  '

# Row 38
$ws.Range("E38").Value = '  -0.02%  '

# Row 39
$ws.Range("B39").Value = 'PolygonEcosystemToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.395'
$ws.Range("E39").Value = '  +4.14%  '

# Row 40
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '163.37'
$ws.Range("E40").Value = '  -0.37%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '195.39'
$ws.Range("E41").Value = '  +7.14%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '20.06'
$ws.Range("E42").Value = '  +0.45%  '

# Row 43
$ws.Range("B43").Value = 'Cronos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.108'
$ws.Range("E43").Value = '  +3.20%  '

# Row 44
$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  +0.09%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.35'
$ws.Range("E45").Value = '  +9.15%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.798'
$ws.Range("E46").Value = '  +13.49%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.75'
$ws.Range("E47").Value = '  +7.47%  '

# Row 48
$ws.Range("B48").Value = 'ImmutableX'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.28'
$ws.Range("E48").Value = '  +5.07%  '

# Row 49
$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '42.71'
$ws.Range("E49").Value = '  -0.78%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.53'
$ws.Range("E50").Value = '  +10.37%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.618'
$ws.Range("E51").Value = '  +6.24%  '
